$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the student email addresses in C2:C4 with a generic placeholder
# email (the underlying mailto: hyperlink targets are left as-is).
$newEmail = "none.none@georgiancollege.ca"

$ws.Range("C2").Value = $newEmail
$ws.Range("C3").Value = $newEmail
$ws.Range("C4").Value = $newEmail
